$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "dimanche, mars 26, 2023"
$ws.Range("B19").Value = "Brad Peat"
$ws.Range("C19").Value = 9

$ws.Range("A20").Value = "samedi, avril 01, 2023"
$ws.Range("B20").Value = "EQUILUX II"
$ws.Range("C20").Value = 20

$ws.Range("E18").Select()
